# Daily attendance processing - 2026-01-09 11:05:33
# For every row in the "Recorded By" column (G), when the cell contains a
# comma-separated list of more than one recorder, rotate the list left by
# one position (move the first entry to the end of the list). Single-entry
# cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }

    $parts = $current -split ",\s*"

    if ($parts.Count -gt 1) {
        $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
        $cell.Value = [string]::Join(", ", $rotated)
    }
}
